$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.580.40'
$ws.Range('E2').Value = '  -1.11%  '
$ws.Range('D3').Value = '1.595.12'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').Value = '  +0.53%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.14'
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.503'
$ws.Range('E6').Value = '  -3.64%  '
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.27'
$ws.Range('E8').Value = '  -4.83%  '
$ws.Range('E9').Value = '  -2.05%  '
$ws.Range('E11').Value = '  -1.91%  '
$ws.Range('D12').Value = '1.822.74'
$ws.Range('E12').Value = '  -2.08%  '
$ws.Range('D13').Value = '1.592.46'
$ws.Range('E13').Value = '  -2.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.86'
$ws.Range('E14').Value = '  -4.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.535'
$ws.Range('E15').Value = '  -4.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.37'
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('D17').Value = '27.568.78'
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '217.00'
$ws.Range('E18').Value = '  -5.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.35'
$ws.Range('E19').Value = '  -4.06%  '
$ws.Range('D20').Value = '0.0₃0693'
$ws.Range('E20').Value = '  -3.77%  '
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.17'
$ws.Range('E22').Value = '  -3.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.64'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.00'
$ws.Range('E24').Value = '  -2.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.96'
$ws.Range('E25').Value = '  -0.71%  '
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.06'
$ws.Range('E28').Value = '  -2.95%  '
$ws.Range('E29').Value = '  -4.07%  '
$ws.Range('E30').Value = '  -1.35%  '
$ws.Range('E31').Value = '  -3.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.25'
$ws.Range('E32').Value = '  -4.98%  '
$ws.Range('D33').Value = '1.368.05'
$ws.Range('E33').Value = '  -1.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.95'
$ws.Range('E34').Value = '  -5.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.52'
$ws.Range('E35').Value = '  -3.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.973'
$ws.Range('E36').Value = '  -3.02%  '
$ws.Range('E37').Value = '  -1.03%  '
$ws.Range('E38').Value = '  -3.19%  '
$ws.Range('E39').Value = '  -3.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.811'
$ws.Range('E40').Value = '  -4.68%  '
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.974'
$ws.Range('E42').Value = '  -3.76%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.79'
$ws.Range('E43').Value = '  -2.51%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.36'
$ws.Range('E44').Value = '  -1.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.98'
$ws.Range('E45').Value = '  -2.70%  '
$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.16'
$ws.Range('E46').Value = '  +0.69%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.733.44'
$ws.Range('E47').Value = '  -2.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.50'
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('E49').Value = '  -3.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0969'
$ws.Range('E50').Value = '  -4.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0497'
$ws.Range('E51').Value = '  -1.28%  '
